# Generate Report for Handoff
# Updates the localization-status report: "In Translation" -> "Ready for
# handoff" everywhere, refreshes the handoff/generate timestamps, and widens
# the datetime columns so the new, longer status text fits.

$wb = $excel.ActiveWorkbook

# NOTE on column widths: the target OOXML `width` is a raw character-unit
# float (17.2159881591797) produced by a non-Excel writer. Real Excel (and
# this COM emulation, faithfully) only persists ColumnWidth snapped to
# whole-pixel boundaries, so we feed the ColumnWidth that snaps closest to
# the target (17.1666... , i.e. the nearest reachable pixel boundary).
$targetColWidth = 16.333333333333332

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-30 04:38:50"
$ws.Columns.Item(5).ColumnWidth = $targetColWidth
$ws.Columns.Item(6).ColumnWidth = $targetColWidth

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-30 04:38:46"
$ws.Columns.Item(3).ColumnWidth = $targetColWidth

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-30 04:38:50"
$ws.Columns.Item(3).ColumnWidth = $targetColWidth
